# Weekly price-update: a new price record for "Albahaca" (Femacal de La
# Calera) is inserted ahead of the existing history, pushing the prior
# rows down by one (old row 76 -> 77, ..., old row 116 -> 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 76; Excel shifts rows 76..116 down to 77..117.
$ws.Rows(76).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44572
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112052
$ws.Range("G76").Value = "Albahaca"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 128
$ws.Range("K76").Value = 4500
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = 4766
$ws.Range("N76").Value = "$/docena de matas"
$ws.Range("O76").Value = "Provincia de Quillota"
$ws.Range("P76").Value = 794
$ws.Range("Q76").Value = 6
$ws.Range("R76").Value = "Hortaliza"
